# Sanity Test - Profile Details
# Applies the changes described by the commit: adds an "Availability"
# column to the EditShareSkill and SignIn sheets, and appends two new
# worksheets ("ProfileDetails" and "SearchSkill").

$wb = $excel.ActiveWorkbook

$editShareSkill = $wb.Worksheets.Item("EditShareSkill")
$signIn         = $wb.Worksheets.Item("SignIn")

# ---------------------------------------------------------------
# Create the two new sheets up front (empty) so every sheet we will
# touch already exists; the actual cell values are filled in below
# in a specific order so the shared-string table is built up the
# same way the original workbook's was.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$profileDetails = $wb.Worksheets.Add($null, $lastSheet)
$profileDetails.Name = "ProfileDetails"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$searchSkill = $wb.Worksheets.Add($null, $lastSheet2)
$searchSkill.Name = "SearchSkill"

# ---------------------------------------------------------------
# Fill in cell values (order chosen so new shared strings are
# interned in the same sequence as the target workbook).
# ---------------------------------------------------------------
$editShareSkill.Range("Q1").Value = 'Availability'
$profileDetails.Range("B1").Value = 'Hours'
$profileDetails.Range("B2").Value = 'As needed'
$profileDetails.Range("C2").Value = 'Less than $500 per month'
$profileDetails.Range("C1").Value = 'EarnTarget'
$editShareSkill.Range("Q2").Value = 'Part Time'
$profileDetails.Range("A2").Value = 'Full Time'
$searchSkill.Range("A1").Value = 'Search'
$searchSkill.Range("A2").Value = 'Java'

$profileDetails.Range("A1").Value = 'Availability'
$signIn.Range("D1").Value = 'Availability'
$signIn.Range("D2").Value = 'Part Time'

# ---------------------------------------------------------------
# Formatting: D1 on SignIn picks up the same header formatting as
# C1 (blue fill header style).
# ---------------------------------------------------------------
$signIn.Range("C1").Copy()
$signIn.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Column widths for the new sheets.
# ---------------------------------------------------------------
$profileDetails.Columns.Item(1).ColumnWidth = 10.619791666666666
$profileDetails.Columns.Item(2).ColumnWidth = 11.709635416666666
$profileDetails.Columns.Item(3).ColumnWidth = 21.893229166666668

$searchSkill.Columns.Item(1).ColumnWidth = 11.529947916666666
$searchSkill.Columns.Item(2).ColumnWidth = 13.799479166666666
$searchSkill.Columns.Item(3).ColumnWidth = 25.072916666666668

# ---------------------------------------------------------------
# Selections / active sheet (mirrors the final view state saved in
# the workbook).
# ---------------------------------------------------------------
$editShareSkill.Activate()
$editShareSkill.Range("Q2").Select()

$signIn.Activate()
$signIn.Range("D2").Select()

$profileDetails.Activate()
$profileDetails.Range("A2").Select()

$searchSkill.Activate()
$searchSkill.Range("A2").Select()
